$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Luca Stoppi results update (First Game results - official) ---
$ws.Range("J2").Value = 2
$ws.Range("N2").Value = 10.96
$ws.Range("O2").Value = 263

# --- Rows 24-26: re-ordered / updated standings ---
# Row 24 -> Coyote (Cri)
$ws.Range("A24").Value = "Coyote (Cri)"
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = 10
$ws.Range("D24").Value = 2
$ws.Range("E24").Value = "'20%"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 93
$ws.Range("J24").Value = 1
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0.5
$ws.Range("M24").Value = 5
$ws.Range("N24").Value = 6.8
$ws.Range("O24").Value = 68
$ws.Range("P24").Value = -4

# Row 25 -> Emilano (Dani)
$ws.Range("A25").Value = "Emilano (Dani)"
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = "'29%"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = 57
$ws.Range("J25").Value = 1
$ws.Range("K25").Value = 1
$ws.Range("L25").Value = 1.86
$ws.Range("M25").Value = 13
$ws.Range("N25").Value = 9.57
$ws.Range("O25").Value = 67
$ws.Range("P25").Value = 0

# Row 26 -> Cerro
$ws.Range("A26").Value = "Cerro"
$ws.Range("B26").Value = 24
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = "'100%"
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 17
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 1
$ws.Range("M26").Value = 5
$ws.Range("N26").Value = 13.4
$ws.Range("O26").Value = 67
$ws.Range("P26").Value = -7
